$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = "S51"
$ws.Range("E8").Value = "Q476"
$ws.Range("D9").Value = "S51"
$ws.Range("E9").Value = "Q476"
$ws.Range("H10").Value = "S51"
$ws.Range("I10").Value = "Q476"
$ws.Range("D11").Value = "S51"
$ws.Range("E11").Value = "Q476"
$ws.Range("D12").Value = "S51"
$ws.Range("E12").Value = "Q476"
$ws.Range("D13").Value = "S51"
$ws.Range("E13").Value = "Q476"
$ws.Range("D14").Value = "S51"
$ws.Range("E14").Value = "Q476"
$ws.Range("D23").Value = "S51"
$ws.Range("E23").Value = "Q475"
$ws.Range("D24").Value = "S51"
$ws.Range("E24").Value = "Q475"
$ws.Range("H25").Value = "S51"
$ws.Range("I25").Value = "Q475"
$ws.Range("D26").Value = "S51"
$ws.Range("E26").Value = "Q475"
$ws.Range("D27").Value = "S51"
$ws.Range("E27").Value = "Q475"
$ws.Range("D28").Value = "S51"
$ws.Range("E28").Value = "Q475"
$ws.Range("D29").Value = "S51"
$ws.Range("E29").Value = "Q475"
$ws.Range("D38").Value = "S51"
$ws.Range("E38").Value = "Q474"
$ws.Range("D39").Value = "S51"
$ws.Range("E39").Value = "Q474"
$ws.Range("D40").Value = "S51"
$ws.Range("E40").Value = "Q474"
$ws.Range("D41").Value = "S51"
$ws.Range("E41").Value = "Q474"
$ws.Range("H42").Value = "S51"
$ws.Range("I42").Value = "Q474"
$ws.Range("D43").Value = "S51"
$ws.Range("E43").Value = "Q474"
$ws.Range("D44").Value = "S51"
$ws.Range("E44").Value = "Q474"
$ws.Range("D45").Value = "S51"
$ws.Range("E45").Value = "Q474"
$ws.Range("D46").Value = "S51"
$ws.Range("E46").Value = "Q474"
$ws.Range("D47").Value = "S51"
$ws.Range("E47").Value = "Q474"
$ws.Range("D56").Value = "S51"
$ws.Range("E56").Value = "Q400"
$ws.Range("D57").Value = "S51"
$ws.Range("E57").Value = "Q400"
$ws.Range("H58").Value = "S51"
$ws.Range("I58").Value = "Q400"
$ws.Range("D59").Value = "S51"
$ws.Range("E59").Value = "Q400"
$ws.Range("D60").Value = "S51"
$ws.Range("E60").Value = "Q400"
$ws.Range("D61").Value = "S51"
$ws.Range("E61").Value = "Q400"
$ws.Range("D62").Value = "S51"
$ws.Range("E62").Value = "Q400"
$ws.Range("D63").Value = "S51"
$ws.Range("E63").Value = "Q400"
$ws.Range("D74").Value = "S51"
$ws.Range("E74").Value = "Q381"
$ws.Range("D75").Value = "S51"
$ws.Range("E75").Value = "Q381"
$ws.Range("D76").Value = "S51"
$ws.Range("E76").Value = "Q381"
$ws.Range("D77").Value = "S51"
$ws.Range("E77").Value = "Q381"
$ws.Range("D78").Value = "S51"
$ws.Range("E78").Value = "Q381"
$ws.Range("H79").Value = "S51"
$ws.Range("I79").Value = "Q381"
$ws.Range("D80").Value = "S51"
$ws.Range("E80").Value = "Q381"
$ws.Range("D81").Value = "S51"
$ws.Range("E81").Value = "Q381"
$ws.Range("D82").Value = "S51"
$ws.Range("E82").Value = "Q381"
$ws.Range("D83").Value = "S51"
$ws.Range("E83").Value = "Q381"
$ws.Range("D84").Value = "S51"
$ws.Range("E84").Value = "Q381"
$ws.Range("D85").Value = "S51"
$ws.Range("E85").Value = "Q381"
$ws.Range("D86").Value = "S51"
$ws.Range("E86").Value = "Q381"
$ws.Range("D87").Value = "S51"
$ws.Range("E87").Value = "Q381"
$ws.Range("D88").Value = "S51"
$ws.Range("E88").Value = "Q381"
$ws.Range("D89").Value = "S51"
$ws.Range("E89").Value = "Q381"
$ws.Range("D90").Value = "S51"
$ws.Range("E90").Value = "Q381"
$ws.Range("D100").Value = "S51"
$ws.Range("E100").Value = "Q470"
$ws.Range("D101").Value = "S51"
$ws.Range("E101").Value = "Q470"
$ws.Range("D102").Value = "S51"
$ws.Range("E102").Value = "Q470"
$ws.Range("D103").Value = "S51"
$ws.Range("E103").Value = "Q470"
$ws.Range("D104").Value = "S51"
$ws.Range("E104").Value = "Q470"
$ws.Range("D105").Value = "S51"
$ws.Range("E105").Value = "Q470"
$ws.Range("H106").Value = "S51"
$ws.Range("I106").Value = "Q470"
$ws.Range("D108").Value = "S51"
$ws.Range("E108").Value = "Q470"
$ws.Range("D109").Value = "S51"
$ws.Range("E109").Value = "Q470"
$ws.Range("D110").Value = "S51"
$ws.Range("E110").Value = "Q470"
$ws.Range("D111").Value = "S51"
$ws.Range("E111").Value = "Q470"
$ws.Range("D112").Value = "S51"
$ws.Range("E112").Value = "Q470"
$ws.Range("D113").Value = "S51"
$ws.Range("E113").Value = "Q470"
$ws.Range("D114").Value = "S51"
$ws.Range("E114").Value = "Q470"
$ws.Range("D115").Value = "S51"
$ws.Range("E115").Value = "Q470"
